# NIT-9019005778.xlsx — "Estado de cuenta" update
#
# Adds a new debt-period row (period 2508) to the worker table and
# re-sorts the existing four periods (2507/2506/2505/2504) into
# ascending order (2504..2507), updating the two summary fields that
# depend on the period count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at 20. This pushes the trailing
#    "divider + signature" rows (previously 24/25) down to 25/26 and
#    keeps row 19's original content/format intact for the moment.
$ws.Rows("20:20").Insert()

# 2) The just-inserted row 20 is blank; row 19 still holds the old
#    "last row" look (thicker bottom border) and its 2504 data.
#    Duplicate row 19 (content + formatting) down into row 20 so the
#    new last row gets the correct border treatment.
$ws.Range("B19:J19").Copy($ws.Range("B20:J20"))

# 3) Row 19 should now look like a normal interior row (same style as
#    rows 16-18), so copy just the formatting from row 16 onto it.
$ws.Range("B16:J16").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)  # xlPasteFormats

# 4) Re-enter the period values so the table reads in ascending order
#    top to bottom, finishing with the brand-new period 2508.
$ws.Range("E16").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"
$ws.Range("E20").Value = "2508"

# 5) Update the summary fields: total overdue value and period count
#    (5 periods x $56,940 = $284,700).
$ws.Range("E11").Value = 284700
$ws.Range("F13").Value = 5
